$wb = $excel.ActiveWorkbook

# Sheet "max-arrecad": rotate categoria_mencao values in A16:A19
$wsMax = $wb.Worksheets.Item("max-arrecad")
$wsMax.Range("A16").Value = "ccxp"
$wsMax.Range("A17").Value = "hqmix"
$wsMax.Range("A18").Value = "questoes_genero"
$wsMax.Range("A19").Value = "angelo_agostini"

# Sheet "tx-sucesso": swap categoria_mencao values in A5/A6 and A9/A10
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A5").Value = "questoes_genero"
$wsTx.Range("A6").Value = "saloes_humor"
$wsTx.Range("A9").Value = "zine"
$wsTx.Range("A10").Value = "lgbtqiamais"
